# feat: add 2022-Q1 data
#
# The old "总计" (summary) sheet becomes the new "2022-Q1" sheet (holding the
# per-fund holdings detail for 2022-Q1), and a brand new "总计" sheet is
# appended after it, carrying the updated roll-up table (with the 2022-Q1
# row inserted at the top).
#
# A tiny scratch sheet is used as scratch space so that numeric-looking
# text (fund codes like "002685", or formatted numbers like "59.37") can be
# written as genuine text cells (matching the source data) instead of being
# auto-coerced to numbers by `Range.Value`.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the current "总计" sheet to "2022-Q1", then add a fresh "总计"
#    sheet right after it, and a throw-away scratch sheet after that.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("总计").Name = "2022-Q1"
$wb.Worksheets.Add($null, $wb.Worksheets.Item("2022-Q1")).Name = "总计"
$wb.Worksheets.Add($null, $wb.Worksheets.Item("总计")).Name = "__scratch__"

# ---------------------------------------------------------------------
# Helper: write a value into a sheet as a plain text cell (no quote-prefix,
# no leftover style) by staging it on the scratch sheet and copying only
# the *value* across.
# ---------------------------------------------------------------------
function Set-TextCell {
    param($sheetName, $row, $col)

    $helper = $wb.Worksheets.Item("__scratch__").Range("A1")
    $helper.Value = "'" + $script:cellText
    $helper.Copy()
    $wb.Worksheets.Item($sheetName).Cells.Item($row, $col).PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 2. Rebuild "2022-Q1" (the per-fund holdings detail for the new quarter).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# columns: code, name, scale, position, ratio, marketvalue, rank
$fundRows = @(
    @("002685", "中欧丰泓沪港深灵活配置混合A",      "59.37", "94.70", "7.79", "4.6249", 3),
    @("013991", "中欧港股通精选一年持有混合A",      "12.87", "94.50", "6.43", "0.8275", 2),
    @("002686", "中欧丰泓沪港深灵活配置混合C",       "7.65", "94.70", "7.79", "0.5959", 3),
    @("013992", "中欧港股通精选一年持有混合C",       "5.32", "94.50", "6.43", "0.3421", 2),
    @("002662", "前海开源沪港深大消费主题混合A",      "1.03", "86.97", "5.92", "0.0610", 8),
    @("007132", "长城港股通价值精选多策略混合",       "0.96", "85.80", "2.52", "0.0242", 10),
    @("002663", "前海开源沪港深大消费主题混合C",      "0.39", "86.97", "5.92", "0.0231", 8),
    @("005255", "浦银安盛港股通量化优选灵活配置混合", "0.35", "90.53", "3.87", "0.0135", 10)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q1.Cells.Item($r, 1).Value = $i

    $script:cellText = $row[0]
    Set-TextCell "2022-Q1" $r 2

    $q1.Cells.Item($r, 3).Value = $row[1]

    $script:cellText = $row[2]
    Set-TextCell "2022-Q1" $r 4

    $script:cellText = $row[3]
    Set-TextCell "2022-Q1" $r 5

    $script:cellText = $row[4]
    Set-TextCell "2022-Q1" $r 6

    $script:cellText = $row[5]
    Set-TextCell "2022-Q1" $r 7

    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 3. Rebuild "总计" (roll-up table): same as before, with a new 2022-Q1
#    row inserted right after the header.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Cells.Clear()

$tot.Range("B1").Value = "日期"
$tot.Range("C1").Value = "持有数量(只)"
$tot.Range("D1").Value = "持有市值(亿元)"

$totRows = @(
    @("2022-Q1", 8, 6.51),
    @("2021-Q4", 4, 1.1),
    @("2021-Q3", 5, 4.31),
    @("2021-Q2", 2, 2.91),
    @("2021-Q1", 2, 1.32),
    @("2020-Q4", 2, 2.32)
)

for ($i = 0; $i -lt $totRows.Length; $i++) {
    $r = $i + 2
    $row = $totRows[$i]
    $tot.Cells.Item($r, 1).Value = $i
    $tot.Cells.Item($r, 2).Value = $row[0]
    $tot.Cells.Item($r, 3).Value = $row[1]
    $tot.Cells.Item($r, 4).Value = $row[2]
}

# ---------------------------------------------------------------------
# 4. Re-apply the standard header / first-column formatting (style index 2
#    in the original workbook: bold, centered, thin-bordered) by copying it
#    from a sheet that already uses it, since `Cells.Clear()` above wiped
#    formatting along with content.
# ---------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item("2021-Q1")

$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

$styleSrc.Range("B1:D1").Copy()
$tot.Range("B1:D1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$tot.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Drop the scratch sheet used for text staging.
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("__scratch__").Delete()
